$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of key terms to append to the "Key Terms" list.
# Each entry: Term (column A), Definition (column B)
$terms = @(
    @("Context Overload", "Providing too much context and instructions can overwhelm the AI, resulting in poor or incomplete responses"),
    @("Task Breakdown", "Breaking a complex request into smaller, more manageable prompts for AI. This avoid overloading it."),
    @("Chain of Thought", "Prompting the AI to explain its reasoning step-by-step. Useful for getting more detailed responses"),
    @("Limitations", "AI assistants have constraints on knowledge, skills, and processing ability. Pushing too hard leads to poor results"),
    @("Rephrasing", "Restating or rewording a prompt to get better results from the AI. Helps clarify intent.")
)

$startRow = 13

for ($i = 0; $i -lt $terms.Length; $i++) {
    $row = $startRow + $i
    $term = $terms[$i][0]
    $def = $terms[$i][1]

    $cellA = $ws.Cells.Item($row, 1)
    $cellB = $ws.Cells.Item($row, 2)

    $cellA.Value = $term
    $cellB.Value = $def

    # Alternate fill color between yellow and green for column A, as in the
    # existing rows above, while column B stays unfilled.
    if (($i % 2) -eq 0) {
        $cellA.Interior.Color = 65535
    } else {
        $cellA.Interior.Color = 5296274
    }
    $cellA.Font.Bold = $true

    # Left/right thin borders only (no top/bottom) for both columns.
    $cellA.Borders.Item(7).LineStyle = 1
    $cellA.Borders.Item(7).Weight = 2
    $cellA.Borders.Item(10).LineStyle = 1
    $cellA.Borders.Item(10).Weight = 2

    $cellB.Borders.Item(7).LineStyle = 1
    $cellB.Borders.Item(7).Weight = 2
    $cellB.Borders.Item(10).LineStyle = 1
    $cellB.Borders.Item(10).Weight = 2
}

$ws.Range("C17").Select()
